$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '69.913.37'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.541.41'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '605.36'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '194.79'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -3.54%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.47'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000303'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.51'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').Value = '4.100.01'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '592.45'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.17'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '12.80'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = '70.013.72'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = '3.528.44'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.986'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '103.15'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.13'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.06'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.25%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.75'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.89%  '
$ws.Range('E28').Value = '  -2.94%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '33.18'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.06%  '
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.23'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -6.39%  '
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.115'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '63.42'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.23'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +6.67%  '
$ws.Range('D36').Value = '3.838.82'
$ws.Range('D37').Value = '0.0₃0823'
$ws.Range('E37').Value = '  +3.80%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '512.55'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '36.50'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('E43').Value = '  -2.84%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0449'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.82'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.139'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.31'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.000245'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.99%  '
$ws.Range('E51').Value = '  +1.85%  '
